$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated l-mixed level values in column J (J4:J43)
$ws.Range("J4").Value = 0.93506
$ws.Range("J5").Value = 0.89353000000000005
$ws.Range("J6").Value = 0.92720999999999998
$ws.Range("J7").Value = 0.89563000000000004
$ws.Range("J8").Value = 0.92318999999999996
$ws.Range("J9").Value = 0.93901999999999997
$ws.Range("J10").Value = 0.87756000000000001
$ws.Range("J11").Value = 0.93145999999999995
$ws.Range("J12").Value = 0.90483000000000002
$ws.Range("J13").Value = 0.92825000000000002
$ws.Range("J14").Value = 0.86412999999999995
$ws.Range("J15").Value = 0.89712000000000003
$ws.Range("J16").Value = 0.91886999999999996
$ws.Range("J17").Value = 0.85948000000000002
$ws.Range("J18").Value = 0.81681999999999999
$ws.Range("J19").Value = 0.90776000000000001
$ws.Range("J20").Value = 0.91281999999999996
$ws.Range("J21").Value = 0.93188000000000004
$ws.Range("J22").Value = 0.43247000000000002
$ws.Range("J23").Value = 0.89351999999999998
$ws.Range("J24").Value = 0.92318
$ws.Range("J25").Value = 0.93506
$ws.Range("J26").Value = 0.83423999999999998
$ws.Range("J27").Value = 0.78149999999999997
$ws.Range("J28").Value = 0.93903999999999999
$ws.Range("J29").Value = 0.87756000000000001
$ws.Range("J30").Value = 0.93186999999999998
$ws.Range("J31").Value = 0.93147000000000002
$ws.Range("J32").Value = 0.89563000000000004
$ws.Range("J33").Value = 0.90483999999999998
$ws.Range("J34").Value = 0.89563999999999999
$ws.Range("J35").Value = 0.78283999999999998
$ws.Range("J36").Value = 0.71048
$ws.Range("J37").Value = 0.86414000000000002
$ws.Range("J38").Value = 0.93506999999999996
$ws.Range("J39").Value = 0.87392999999999998
$ws.Range("J40").Value = 0.92986000000000002
$ws.Range("J41").Value = 0.88553000000000004
$ws.Range("J42").Value = 0.92818999999999996
$ws.Range("J43").Value = 0.91883000000000004

# Selection changed to whole column B (B1:B1048576) with active cell B1
$null = $ws.Columns("B:B").Select()

# Column I (hidden) width restored to sheet default width (best achievable via ColumnWidth)
$ws.Columns.Item(9).ColumnWidth = 10
